$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated save_data: the "K" column (G) now reports strikeouts-derived
# values computed from the new S0/SF logic ("use K instead of Strike#")
# instead of the old Strike# count. Write the recalculated values per row.
$newK = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 2
    6  = 1
    7  = 2
    8  = 2
    9  = 1
    10 = 0
    11 = 1
    13 = 1
    14 = 0
    15 = 2
    16 = 0
    17 = 2
    18 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 0
    24 = 1
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 2
    30 = 1
    31 = 2
    32 = 2
    33 = 1
    34 = 0
    35 = 0
    36 = 2
    38 = 2
    39 = 1
    40 = 2
    41 = 3
    42 = 3
    43 = 1
    44 = 0
    45 = 1
    46 = 1
    47 = 2
    48 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
